$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue "D2" "24.600.17"
Set-TextValue "E2" "  -1.01%  "
Set-TextValue "D3" "1.651.99"
Set-TextValue "E3" "  -3.08%  "
Set-TextValue "D4" "1.004"
Set-TextValue "E4" "  -0.23%  "
Set-TextValue "D5" "318.01"
Set-TextValue "E5" "  +1.90%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.11%  "
Set-TextValue "D7" "0.3628"
Set-TextValue "E7" "  -2.71%  "
Set-TextValue "D8" "46.60"
Set-TextValue "E8" "  -5.46%  "
Set-TextValue "D9" "0.3249"
Set-TextValue "D10" "1.129"
Set-TextValue "E10" "  -7.26%  "
Set-TextValue "D11" "0.07025"
Set-TextValue "D12" "1.001"
Set-TextValue "E12" "  -0.12%  "
Set-TextValue "D13" "5.999"
Set-TextValue "E13" "  -5.04%  "
Set-TextValue "D14" "19.45"
Set-TextValue "E14" "  -7.72%  "
Set-TextValue "D15" "1.654.10"
Set-TextValue "D16" "6.590"
Set-TextValue "E16" "  -6.28%  "
Set-TextValue "E17" "  -8.02%  "
Set-TextValue "D18" "0.06613"
Set-TextValue "E18" "  -1.58%  "
Set-TextValue "E19" "  +0.20%  "
Set-TextValue "D20" "78.58"
Set-TextValue "E20" "  -6.27%  "
Set-TextValue "D21" "5.898"
Set-TextValue "E21" "  -7.14%  "
Set-TextValue "D22" "15.67"
Set-TextValue "E22" "  -9.29%  "
Set-TextValue "D23" "12.52"
Set-TextValue "E23" "  -4.44%  "
Set-TextValue "D24" "24.568.48"
Set-TextValue "E24" "  -1.20%  "
Set-TextValue "D25" "2.436"
Set-TextValue "E25" "  -0.86%  "
Set-TextValue "D26" "2.367"
Set-TextValue "E26" "  -14.73%  "
Set-TextValue "D27" "147.73"
Set-TextValue "E27" "  -1.17%  "
Set-TextValue "E28" "  -8.79%  "
Set-TextValue "D29" "1.835.19"
Set-TextValue "E29" "  -3.15%  "
Set-TextValue "D30" "1.202"
Set-TextValue "E30" "  -4.20%  "
Set-TextValue "D31" "125.17"
Set-TextValue "E31" "  -5.49%  "
Set-TextValue "D32" "4.067"
Set-TextValue "E32" "  -3.86%  "
Set-TextValue "D33" "5.782"
Set-TextValue "E33" "  -14.61%  "
Set-TextValue "D34" "0.08440"
Set-TextValue "E34" "  -3.69%  "
Set-TextValue "D35" "1.672"
Set-TextValue "E35" "  -5.57%  "
Set-TextValue "D36" "12.19"
Set-TextValue "E36" "  -11.21%  "
Set-TextValue "E37" "  -0.47%  "
Set-TextValue "D38" "5.191"
Set-TextValue "E38" "  -7.05%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.02225"
Set-TextValue "E39" "  -7.63%  "
Set-TextValue "B40" "Hedera"
Set-TextValue "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.05995"
Set-TextValue "E40" "  -9.82%  "
Set-TextValue "D41" "0.2067"
Set-TextValue "E41" "  -7.66%  "
Set-TextValue "D42" "8.153"
Set-TextValue "E42" "  -10.59%  "
Set-TextValue "D43" "1.001"
Set-TextValue "E43" "  +0.02%  "
Set-TextValue "D44" "0.5897"
Set-TextValue "E44" "  -8.19%  "
Set-TextValue "D45" "3.827"
Set-TextValue "E45" "  -0.13%  "
Set-TextValue "D46" "12.66"
Set-TextValue "E46" "  -8.25%  "
Set-TextValue "D47" "0.5604"
Set-TextValue "E47" "  -8.59%  "
Set-TextValue "D48" "123.95"
Set-TextValue "E48" "  -4.47%  "
Set-TextValue "D49" "1.943"
Set-TextValue "E49" "  -7.96%  "
Set-TextValue "D50" "0.06934"
Set-TextValue "E50" "  -5.14%  "
Set-TextValue "D51" "1.189"
Set-TextValue "E51" "  -2.67%  "
